$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 269 (Femacal de La Calera /
# Papa "1a (cosecha)" row dated 2021-06-10), pushing the old rows 269:291
# down to 271:293.
$ws.Rows("269:270").Insert()

# Populate the two freshly inserted rows with the new weekly data
# (date 44461 = 2021-09-22), matching the columns used by every other
# row in this subset: Mercado ID, Mercado, Región, Fecha, Codreg,
# Categoría ID, Categoría, Variedad, Calidad, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Unidad de comercialización,
# Origen, Precio $/Kg, Kg o Unidades, Clasificación.

$ws.Range("A269").Value = 3
$ws.Range("B269").Value = "Femacal de La Calera"
$ws.Range("C269").Value = "Coquimbo"
$ws.Range("D269").Value = 44461
$ws.Range("E269").Value = 5
$ws.Range("F269").Value = 100114001
$ws.Range("G269").Value = "Papa"
$ws.Range("H269").Value = "Asterix"
$ws.Range("I269").Value = "1a (guarda)"
$ws.Range("J269").Value = 250
$ws.Range("K269").Value = 9500
$ws.Range("L269").Value = 10000
$ws.Range("M269").Value = 9760
$ws.Range("N269").Value = "`$/saco 25 kilos"
$ws.Range("O269").Value = "Región del Maule"
$ws.Range("P269").Value = 390
$ws.Range("Q269").Value = 25
$ws.Range("R269").Value = "Hortaliza"

$ws.Range("A270").Value = 3
$ws.Range("B270").Value = "Femacal de La Calera"
$ws.Range("C270").Value = "Coquimbo"
$ws.Range("D270").Value = 44461
$ws.Range("E270").Value = 5
$ws.Range("F270").Value = 100114001
$ws.Range("G270").Value = "Papa"
$ws.Range("H270").Value = "Rosara"
$ws.Range("I270").Value = "1a (guarda)"
$ws.Range("J270").Value = 390
$ws.Range("K270").Value = 9000
$ws.Range("L270").Value = 9500
$ws.Range("M270").Value = 9295
$ws.Range("N270").Value = "`$/saco 25 kilos"
$ws.Range("O270").Value = "Región del Maule"
$ws.Range("P270").Value = 372
$ws.Range("Q270").Value = 25
$ws.Range("R270").Value = "Hortaliza"
